$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer samples captured and prepended to the existing window.
# This pushes the older samples down and drops the oldest ones that fall
# past the bottom of the tracked range (row 21).
$newRows = @(
    @(-3.012916564941406, 8.089370727539062, -0.1633265316486358),
    @(-3.395848751068115, 8.023316383361816,  0.0382503271102905),
    @(-3.384797096252441, 7.934267520904541,  0.07479587197303771),
    @(-3.632324695587158, 7.965863227844238,  0.0220168232917785)
)

$lastRow = 21
$numNew = $newRows.Count

# Shift the existing data rows (2..lastRow) down by $numNew rows, dropping
# whatever would fall past $lastRow.
for ($r = $lastRow - $numNew; $r -ge 2; $r--) {
    $destRow = $r + $numNew
    $ws.Range("A$r`:C$r").Copy() | Out-Null
    $ws.Range("A$destRow`:C$destRow").PasteSpecial(-4104) | Out-Null
}
$excel.CutCopyMode = $false

# Write the newly captured rows into the freed-up space at the top.
for ($i = 0; $i -lt $numNew; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($row, 3).Value = $newRows[$i][2]
}
